# Append Sunday, Jan 15 departures to the "Main Data" sheet (POZ_Departures).
# Mirrors the existing row layout: NUMBER, DATE, TIME, FLIGHT, TO, SHORT,
# AIRLINE, MODEL, AIRCFAT ID, STATUS, (blank), DIFFERENCE, (blank)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 163 (NUMBER 162)
$ws.Cells.Item(163,1).Value = 162
$ws.Cells.Item(163,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(163,3).Value = "5:50 AM"
$ws.Cells.Item(163,4).Value = "E45069"
$ws.Cells.Item(163,5).Value = "Hurghada"
$ws.Cells.Item(163,6).Value = "(HRG)"
$ws.Cells.Item(163,7).Value = "Enter Air "
$ws.Cells.Item(163,8).Value = "B738"
$ws.Cells.Item(163,9).Value = "(SP-ESF)"
$ws.Cells.Item(163,10).Value = "5:56 AM"
$ws.Cells.Item(163,12).Value = "0 hours, 6 minutes"

# Row 164 (NUMBER 163)
$ws.Cells.Item(164,1).Value = 163
$ws.Cells.Item(164,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(164,3).Value = "6:00 AM"
$ws.Cells.Item(164,4).Value = "FR5073"
$ws.Cells.Item(164,5).Value = "Birmingham"
$ws.Cells.Item(164,6).Value = "(BHX)"
$ws.Cells.Item(164,7).Value = "Ryanair "
$ws.Cells.Item(164,8).Value = "B738"
$ws.Cells.Item(164,9).Value = "(SP-RKR)"
$ws.Cells.Item(164,10).Value = "6:04 AM"
$ws.Cells.Item(164,12).Value = "0 hours, 4 minutes"

# Row 165 (NUMBER 164)
$ws.Cells.Item(165,1).Value = 164
$ws.Cells.Item(165,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(165,3).Value = "6:20 AM"
$ws.Cells.Item(165,4).Value = "FR7906"
$ws.Cells.Item(165,5).Value = "Alicante"
$ws.Cells.Item(165,6).Value = "(ALC)"
$ws.Cells.Item(165,7).Value = "Ryanair "
$ws.Cells.Item(165,8).Value = "B738"
$ws.Cells.Item(165,9).Value = "(SP-RSM)"
$ws.Cells.Item(165,10).Value = "6:20 AM"
$ws.Cells.Item(165,12).Value = "0 hours, 0 minutes"

# Row 166 (NUMBER 165)
$ws.Cells.Item(166,1).Value = 165
$ws.Cells.Item(166,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(166,3).Value = "6:45 AM"
$ws.Cells.Item(166,4).Value = "LH1381"
$ws.Cells.Item(166,5).Value = "Frankfurt"
$ws.Cells.Item(166,6).Value = "(FRA)"
$ws.Cells.Item(166,7).Value = "Lufthansa "
$ws.Cells.Item(166,8).Value = "CRJ9"
$ws.Cells.Item(166,9).Value = "(D-ACNB)"
$ws.Cells.Item(166,10).Value = "6:55 AM"
$ws.Cells.Item(166,12).Value = "0 hours, 10 minutes"

# Row 167 (NUMBER 166)
$ws.Cells.Item(167,1).Value = 166
$ws.Cells.Item(167,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(167,3).Value = "11:55 AM"
$ws.Cells.Item(167,4).Value = "FR1751"
$ws.Cells.Item(167,5).Value = "London"
$ws.Cells.Item(167,6).Value = "(STN)"
$ws.Cells.Item(167,7).Value = "Ryanair "
$ws.Cells.Item(167,8).Value = "B738"
$ws.Cells.Item(167,9).Value = "(SP-RKR)"
$ws.Cells.Item(167,10).Value = "11:55 AM"
$ws.Cells.Item(167,12).Value = "0 hours, 0 minutes"

# Row 168 (NUMBER 167)
$ws.Cells.Item(168,1).Value = 167
$ws.Cells.Item(168,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(168,3).Value = "12:55 PM"
$ws.Cells.Item(168,4).Value = "LH1639"
$ws.Cells.Item(168,5).Value = "Munich"
$ws.Cells.Item(168,6).Value = "(MUC)"
$ws.Cells.Item(168,7).Value = "Lufthansa "
$ws.Cells.Item(168,8).Value = "CRJ9"
$ws.Cells.Item(168,9).Value = "(D-ACNH)"
$ws.Cells.Item(168,10).Value = "12:55 PM"
$ws.Cells.Item(168,12).Value = "0 hours, 0 minutes"

# Row 169 (NUMBER 168)
$ws.Cells.Item(169,1).Value = 168
$ws.Cells.Item(169,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(169,3).Value = "2:20 PM"
$ws.Cells.Item(169,4).Value = "FR1975"
$ws.Cells.Item(169,5).Value = "Dublin"
$ws.Cells.Item(169,6).Value = "(DUB)"
$ws.Cells.Item(169,7).Value = "Ryanair "
$ws.Cells.Item(169,8).Value = "B738"
$ws.Cells.Item(169,9).Value = "(SP-RSM)"
$ws.Cells.Item(169,10).Value = "2:27 PM"
$ws.Cells.Item(169,12).Value = "0 hours, 7 minutes"

# Row 170 (NUMBER 169)
$ws.Cells.Item(170,1).Value = 169
$ws.Cells.Item(170,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(170,3).Value = "2:45 PM"
$ws.Cells.Item(170,4).Value = "LO3944"
$ws.Cells.Item(170,5).Value = "Warsaw"
$ws.Cells.Item(170,6).Value = "(WAW)"
$ws.Cells.Item(170,7).Value = "LOT "
$ws.Cells.Item(170,8).Value = "E170"
$ws.Cells.Item(170,9).Value = "(SP-LDI)"
$ws.Cells.Item(170,10).Value = "3:04 PM"
$ws.Cells.Item(170,12).Value = "0 hours, 19 minutes"

# Row 171 (NUMBER 170)
$ws.Cells.Item(171,1).Value = 170
$ws.Cells.Item(171,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(171,3).Value = "3:05 PM"
$ws.Cells.Item(171,4).Value = "KL1274"
$ws.Cells.Item(171,5).Value = "Amsterdam"
$ws.Cells.Item(171,6).Value = "(AMS)"
$ws.Cells.Item(171,7).Value = "KLM "
$ws.Cells.Item(171,8).Value = "E75L"
$ws.Cells.Item(171,9).Value = "(PH-EXR)"
$ws.Cells.Item(171,10).Value = "3:02 PM"
$ws.Cells.Item(171,12).Value = "0 hours, -3 minutes"

# Row 172 (NUMBER 171)
$ws.Cells.Item(172,1).Value = 171
$ws.Cells.Item(172,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(172,3).Value = "3:25 PM"
$ws.Cells.Item(172,4).Value = "LH1391"
$ws.Cells.Item(172,5).Value = "Frankfurt"
$ws.Cells.Item(172,6).Value = "(FRA)"
$ws.Cells.Item(172,7).Value = "Lufthansa "
$ws.Cells.Item(172,8).Value = "CRJ9"
$ws.Cells.Item(172,9).Value = "(D-ACNW)"
$ws.Cells.Item(172,10).Value = "3:32 PM"
$ws.Cells.Item(172,12).Value = "0 hours, 7 minutes"

# Row 173 (NUMBER 172)
$ws.Cells.Item(173,1).Value = 172
$ws.Cells.Item(173,2).Value = "Sunday, Jan 15"
$ws.Cells.Item(173,3).Value = "3:50 PM"
$ws.Cells.Item(173,4).Value = "UNKNOWN"
$ws.Cells.Item(173,5).Value = "Gothenburg"
$ws.Cells.Item(173,6).Value = "(GOT)"
$ws.Cells.Item(173,7).Value = "Ryanair "
$ws.Cells.Item(173,8).Value = "B738"
$ws.Cells.Item(173,9).Value = "(SP-RSX)"
$ws.Cells.Item(173,10).Value = "3:53 PM"
$ws.Cells.Item(173,12).Value = "0 hours, 3 minutes"
